# Update row 8 (ano 2025) metrics on Sheet1 to reflect refreshed "BIBI" data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value = 1480
$ws.Range("D8").Value = 228
$ws.Range("E8").Value = 1252
$ws.Range("F8").Value = 9.351927809680067
$ws.Range("G8").Value = 84.5945945945946
$ws.Range("H8").Value = 15.40540540540541
